$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-27 Friday", "2024-12-28 Saturday"),
    @("676÷9=", "909÷6="),
    @("892÷6=", "220÷3="),
    @("167÷5=", "568÷6="),
    @("406÷3=", "276÷2="),
    @("106÷9=", "234÷6="),
    @("711÷6=", "288÷8="),
    @("100÷4=", "882÷5="),
    @("410÷3=", "525÷5="),
    @("586÷8=", "121÷2="),
    @("721÷6=", "647÷8="),
    @("677÷3=", "308÷5="),
    @("509÷4=", "565÷2="),
    @("823÷2=", "780÷2="),
    @("237÷3=", "926÷4="),
    @("157÷9=", "515÷3="),
    @("611÷4=", "935÷3="),
    @("812÷8=", "809÷4="),
    @("650÷5=", "278÷3="),
    @("976÷9=", "326÷8="),
    @("577÷9=", "751÷7="),
    @("141÷5=", "468÷6="),
    @("485÷3=", "961÷7="),
    @("188÷9=", "768÷4="),
    @("429÷8=", "915÷9="),
    @("479÷5=", "233÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
